$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cancel Order")

# Clear result values that no longer apply
$ws.Range("D2").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("D8").Value = ""

# Update result values to the new expected API responses
$ws.Range("D4").Value = '"code":401,"message":"Unauthorized"'
$ws.Range("D5").Value = '"code":401,"message":"Unauthorized"'
$ws.Range("D7").Value = '"status":404,"error":"Not Found"'

# Widen column D to fit the new, longer values
# (ColumnWidth uses character units; Excel stores a slightly larger value
# in the file, so compensate to land exactly on the target stored width)
$ws.Columns.Item(4).ColumnWidth = 32.666666666666664

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("D8").Select()
